# Add a new "2023" column (Q) to the table, mirroring the formatting of the
# existing last column (P) and filling in the new data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from column P (rows 3-7) over to column Q,
# so the new cells inherit the same borders/fonts/number formats as the
# rest of the table.
$ws.Range("P3:P7").Copy()
$ws.Range("Q3:Q7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new column's data.
$ws.Range("Q4").Value = 2023
$ws.Range("Q5").Value = 93.7
$ws.Range("Q6").Value = 95.5
$ws.Range("Q7").Value = 97.1

# Restore the selection to the top-left cell (matches a normal save state).
$ws.Range("A1").Select()
